# Generate Report for Handback
#
# For the file "8abf2c0b-ef7a-433f-91fb-842f5c343a46" row (row 7) on both the
# "zh-cn" and "de-de" language sheets, the handback-report generator fills in
# the "Latest Target File" / "Latest Handback File" / "Latest Handback
# DateTime" / "Error Detail" columns (I, J, K, P) because the handback that
# came in was stale (not built from the latest handoff).

$wb = $excel.ActiveWorkbook

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/6f7c5a8fa5b7ddf139f417cb218212b5d457ddb9/e2e/8abf2c0b-ef7a-433f-91fb-842f5c343a46.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2ddcb7e7c4d72d001b24c78db0317f71bb086f39/e2e/8abf2c0b-ef7a-433f-91fb-842f5c343a46.md."

foreach ($sheetName in @("zh-cn", "de-de")) {

    $ws = $wb.Worksheets.Item($sheetName)

    # --- Row 7 values -------------------------------------------------
    # I7: Latest Target File -> the handed-off source markdown file name
    $ws.Range("I7").Value = "8abf2c0b-ef7a-433f-91fb-842f5c343a46.md"

    # J7: Latest Handback File -> same xlf file name already shown in G7
    if ($sheetName -eq "zh-cn") {
        $ws.Range("J7").Value = "8abf2c0b-ef7a-433f-91fb-842f5c343a46.790e001d92fb8e72886fd819eeb473610b46152b.zh-cn.xlf"
        $ws.Range("K7").Value = "2016-09-04 04:47:52"
    } else {
        $ws.Range("J7").Value = "8abf2c0b-ef7a-433f-91fb-842f5c343a46.790e001d92fb8e72886fd819eeb473610b46152b.de-de.xlf"
        $ws.Range("K7").Value = "2016-09-04 04:48:00"
    }

    # P7: Error Detail
    $ws.Range("P7").Value = $errorDetail

    # I7 becomes a hyperlink (same target as A7's handoff link), so rebuild
    # the sheet's hyperlink collection in display order: existing links
    # first, then the new I7 link inserted right after A7's, then the link
    # that used to trail (A8) last - this keeps every relationship id
    # shifting exactly like Excel would when a link is inserted mid-list.
    $links = @()
    foreach ($hl in $ws.Hyperlinks) {
        $links += , @($hl.Range.Address($false, $false), $hl.Address, $hl.TextToDisplay)
    }

    $ws.Hyperlinks.Delete()

    foreach ($link in $links) {
        $cellRef = $link[0]
        $address = $link[1]
        $display = $link[2]

        $ws.Hyperlinks.Add($ws.Range($cellRef), $address, $null, $null, $display)

        if ($cellRef -eq "A7") {
            $ws.Hyperlinks.Add($ws.Range("I7"), $address, $null, $null, "8abf2c0b-ef7a-433f-91fb-842f5c343a46.md")
        }
    }

    # Give I7 the same visual "hyperlink" look (underlined, themed blue font)
    # used by the sheet's other link cells.
    $ws.Range("I7").Font.Underline = 2
    $ws.Range("I7").Font.Color = 15570276

    # --- Column P (16) widens to fit the new Error Detail text --------
    $ws.Columns.Item(16).ColumnWidth = 39.166666666666664
}
